$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.812.05'
$ws.Range("E2").Value = '  -6.22%  '

# Row 3
$ws.Range("D3").Value = '2.983.48'
$ws.Range("E3").Value = '  -6.58%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.93%  '

# Row 8
$ws.Range("D8").Value = '2.976.06'
$ws.Range("E8").Value = '  -6.66%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.65%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.130'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.07%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.06'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.52%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.57%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000218'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.33%  '

# Row 15
$ws.Range("E15").Value = '  +0.80%  '

# Row 16
$ws.Range("D16").Value = '3.475.57'
$ws.Range("E16").Value = '  -6.39%  '

# Row 17
$ws.Range("D17").Value = '2.987.52'
$ws.Range("E17").Value = '  -6.20%  '

# Row 18
$ws.Range("D18").Value = '59.884.61'
$ws.Range("E18").Value = '  -6.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.89%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '425.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.665'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.04%  '

# Row 24
$ws.Range("E24").Value = '  -2.46%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.03%  '

# Row 26
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("E27").Value = '  +0.27%  '

# Row 28
$ws.Range("E28").Value = '  -6.32%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.80%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.78%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.47%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0936'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.42%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.55%  '

# Row 35
$ws.Range("B35").Value = 'Mantle'
$ws.Range("C35").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.927'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.74%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '49.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.52%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -17.20%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0656'
$ws.Range("E38").Value = '  -11.71%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.45%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0354'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.59%  '

# Row 41
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '379.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.12%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.108'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.56%  '

# Row 43
$ws.Range("D43").Value = '2.661.43'
$ws.Range("E43").Value = '  -5.32%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.91%  '

# Row 45
$ws.Range("E45").Value = '  +0.00%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.234'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.26%  '

# Row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.84%  '

# Row 49
$ws.Range("E49").Value = '  -4.34%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.29%  '

# Row 51
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.76%  '
